# Natmi following Dr Hou advice
# Update Ligand-expressing cells (E), Receptor-expressing cells (K) counts,
# and recompute the dependent expression / specificity columns (G,H,M,N,O,P,Q,R,S,T)
# for rows 2-7 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{
    2 = @{ E=3; G=4.619088000000001; H=13.857264; K=2; M=121.2130545; N=242.426109;
           O=0.4742742514031324; P=0.4028886461111009; Q=559.8937654842961; R=3359.362592905776;
           S=0.4742742514031324; T=0.4028886461111009 }
    3 = @{ E=3; G=4.619088000000001; H=13.857264; K=3; M=0.9119573333333332; N=2.735872;
           O=0.003568245048871889; P=0.004546753526507615; Q=4.212411174912; R=37.911700574208;
           S=0.003568245048871889; T=0.004546753526507615 }
    4 = @{ E=3; G=4.619088000000001; H=13.857264; K=3; M=0.7040186666666667; N=2.112056;
           O=0.002754636680714656; P=0.003510031926267592; Q=3.251924174976001; R=29.267317574784;
           S=0.002754636680714656; T=0.003510031926267592 }
    5 = @{ E=3; G=4.619088000000001; H=13.857264; K=3; M=0.6855039999999999; N=2.056512;
           O=0.002682193743693282; P=0.00341772319330189; Q=3.166403300352; R=28.497629703168;
           S=0.002682193743693282; T=0.00341772319330189 }
    6 = @{ E=3; G=4.619088000000001; H=13.857264; K=3; M=88.26666266666666; N=264.799988;
           O=0.3453638350486923; P=0.4400718598158738; Q=407.711482323648; R=3669.403340912832;
           S=0.3453638350486923; T=0.4400718598158738 }
    7 = @{ E=3; G=4.619088000000001; H=13.857264; K=2; M=43.794673; N=87.58934600000001;
           O=0.1713568380748954; P=0.1455649854269483; Q=202.291448518224; R=1213.748691109344;
           S=0.1713568380748954; T=0.1455649854269483 }
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $ws.Range("E$r").Value = $vals.E
    $ws.Range("G$r").Value = $vals.G
    $ws.Range("H$r").Value = $vals.H
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("Q$r").Value = $vals.Q
    $ws.Range("R$r").Value = $vals.R
    $ws.Range("S$r").Value = $vals.S
    $ws.Range("T$r").Value = $vals.T
}
